# daily auto push: 2026-02-10 14:24 UTC
# Insert a new data row at row 805 (pushes the existing 805:846 block down to
# 806:847) and populate it with the day's measurement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 805..846 down by one, creating a blank row 805.
$ws.Rows("805:805").Insert()

# New row's date must stay a literal text string like the rest of column A
# (values are stored as plain "yyyy/mm/dd" text, not real dates), so force
# text entry with a leading apostrophe and then drop the quote-prefix style
# it leaves behind so the cell matches its neighbours (no explicit style).
$ws.Range("A805").Value = "'2026/02/10"
$ws.Range("A805").ClearFormats()

$ws.Range("B805").Value = "火"
$ws.Range("C805").Value = 21
$ws.Range("D805").Value = 201
